# Update "想去人数" (interest count, column F) for a handful of events.
# Each event appears both on its category sheet ("展览") and on the
# aggregate "全部类型" sheet, so the same bump is applied in both places.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet (row => new value)
$wsExhibition.Range("F3").Value  = 1002
$wsExhibition.Range("F5").Value  = 452
$wsExhibition.Range("F6").Value  = 712
$wsExhibition.Range("F10").Value = 399
$wsExhibition.Range("F12").Value = 75
$wsExhibition.Range("F13").Value = 834
$wsExhibition.Range("F15").Value = 1982
$wsExhibition.Range("F16").Value = 476
$wsExhibition.Range("F17").Value = 7091
$wsExhibition.Range("F18").Value = 525
$wsExhibition.Range("F20").Value = 53
$wsExhibition.Range("F21").Value = 90
$wsExhibition.Range("F23").Value = 218

# 全部类型 sheet (same events, different row numbers)
$wsAllTypes.Range("F7").Value  = 1002
$wsAllTypes.Range("F11").Value = 452
$wsAllTypes.Range("F12").Value = 712
$wsAllTypes.Range("F17").Value = 399
$wsAllTypes.Range("F20").Value = 75
$wsAllTypes.Range("F22").Value = 834
$wsAllTypes.Range("F25").Value = 1982
$wsAllTypes.Range("F26").Value = 476
$wsAllTypes.Range("F27").Value = 7091
$wsAllTypes.Range("F29").Value = 525
$wsAllTypes.Range("F31").Value = 53
$wsAllTypes.Range("F32").Value = 90
$wsAllTypes.Range("F35").Value = 218
